$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 837.2909
$ws.Range("J17").Value = 782.1667
$ws.Range("L17").Value = 2346.5001
$ws.Range("N17").Value = -2682.5001
$ws.Range("H28").Value = 652.5417
$ws.Range("I28").Value = 614.5
$ws.Range("J28").Value = 766.6667
$ws.Range("K28").Value = 614.5
$ws.Range("L28").Value = 766.6667
$ws.Range("M28").Value = -129.5
$ws.Range("N28").Value = -1736.6667
$ws.Range("H135").Value = 359.5
$ws.Range("I135").Value = 310.23077
$ws.Range("K135").Value = 2792.07693
$ws.Range("M135").Value = -257.0769300000002
$ws.Range("H138").Value = 3049.3
$ws.Range("I138").Value = 707.13513
$ws.Range("J138").Value = 4424.857
$ws.Range("K138").Value = 2121.40539
$ws.Range("L138").Value = 13274.571
$ws.Range("M138").Value = 3018.59461
$ws.Range("N138").Value = -23554.571
$ws.Range("H140").Value = 50626.668
$ws.Range("J140").Value = 51638.184
$ws.Range("L140").Value = 51638.184
$ws.Range("N140").Value = -61998.184

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 4000
$ws.Range("J4").Value = 4000
$ws.Range("L4").Value = 4000
$ws.Range("N4").Value = -4232
$ws.Range("H32").Value = 5884.9297
$ws.Range("I32").Value = 4629.0244
$ws.Range("J32").Value = 9103.1875
$ws.Range("K32").Value = 4629.0244
$ws.Range("L32").Value = 9103.1875
$ws.Range("M32").Value = -4342.0244
$ws.Range("N32").Value = -9677.1875
$ws.Range("H61").Value = 1202.6666
$ws.Range("I61").Value = 958.2069
$ws.Range("K61").Value = 958.2069
$ws.Range("M61").Value = -746.2069
$ws.Range("H110").Value = 1940.0834
$ws.Range("I110").Value = 1944.375
$ws.Range("J110").Value = 1931.5
$ws.Range("K110").Value = 1944.375
$ws.Range("L110").Value = 1931.5
$ws.Range("M110").Value = 100.625
$ws.Range("N110").Value = -6021.5
$ws.Range("H136").Value = 1202.6666
$ws.Range("I136").Value = 958.2069
$ws.Range("K136").Value = 2874.6207
$ws.Range("M136").Value = -324.6206999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 632.6667
$ws.Range("I5").Value = 632.6667
$ws.Range("K5").Value = 632.6667
$ws.Range("M5").Value = -519.6667
$ws.Range("H105").Value = 1667.4138
$ws.Range("I105").Value = 1629.7307
$ws.Range("K105").Value = 1629.7307
$ws.Range("M105").Value = 117.2692999999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1783.4459
$ws.Range("I58").Value = 1527.1617
$ws.Range("J58").Value = 4688
$ws.Range("K58").Value = 1527.1617
$ws.Range("L58").Value = 4688
$ws.Range("M58").Value = -1324.1617
$ws.Range("N58").Value = -5094
$ws.Range("H99").Value = 9096632
$ws.Range("I99").Value = 18185764
$ws.Range("J99").Value = 7500
$ws.Range("K99").Value = 18185764
$ws.Range("L99").Value = 7500
$ws.Range("M99").Value = -18184266
$ws.Range("N99").Value = -10496
$ws.Range("H126").Value = 9096632
$ws.Range("I126").Value = 18185764
$ws.Range("J126").Value = 7500
$ws.Range("K126").Value = 54557292
$ws.Range("L126").Value = 22500
$ws.Range("M126").Value = -54554822
$ws.Range("N126").Value = -27440
$ws.Range("H132").Value = 2179.0715
$ws.Range("I132").Value = 1110.1
$ws.Range("J132").Value = 4851.5
$ws.Range("K132").Value = 3330.3
$ws.Range("L132").Value = 14554.5
$ws.Range("M132").Value = -800.2999999999997
$ws.Range("N132").Value = -19614.5
$ws.Range("H134").Value = 2968.3684
$ws.Range("I134").Value = 1299.8334
$ws.Range("K134").Value = 3899.5002
$ws.Range("M134").Value = -1364.5002
$ws.Range("H136").Value = 1783.4459
$ws.Range("I136").Value = 1527.1617
$ws.Range("J136").Value = 4688
$ws.Range("K136").Value = 4581.4851
$ws.Range("L136").Value = 14064
$ws.Range("M136").Value = -2031.4851
$ws.Range("N136").Value = -19164
$ws.Range("H138").Value = 46032.223
$ws.Range("J138").Value = 46032.223
$ws.Range("L138").Value = 46032.223
$ws.Range("N138").Value = -56312.223
$ws.Range("H139").Value = 178000
$ws.Range("J139").Value = 178000
$ws.Range("L139").Value = 178000
$ws.Range("N139").Value = -188280
$ws.Range("H140").Value = 77476
$ws.Range("J140").Value = 77476
$ws.Range("L140").Value = 77476
$ws.Range("N140").Value = -87836

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 624.8182
$ws.Range("I18").Value = 196.28572
$ws.Range("J18").Value = 1374.75
$ws.Range("K18").Value = 588.85716
$ws.Range("L18").Value = 4124.25
$ws.Range("M18").Value = -419.85716
$ws.Range("N18").Value = -4462.25
$ws.Range("H131").Value = 6667463
$ws.Range("J131").Value = 862.25757
$ws.Range("L131").Value = 2586.77271
$ws.Range("N131").Value = -12666.77271
$ws.Range("H140").Value = 2412.9614
$ws.Range("I140").Value = 2507.611
$ws.Range("K140").Value = 7522.833
$ws.Range("M140").Value = -2342.833

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 93782.09
$ws.Range("I22").Value = 201840.2
$ws.Range("J22").Value = 3733.6667
$ws.Range("K22").Value = 201840.2
$ws.Range("L22").Value = 3733.6667
$ws.Range("M22").Value = -201545.2
$ws.Range("N22").Value = -4323.6667
$ws.Range("H27").Value = 93782.09
$ws.Range("I27").Value = 201840.2
$ws.Range("J27").Value = 3733.6667
$ws.Range("K27").Value = 201840.2
$ws.Range("L27").Value = 3733.6667
$ws.Range("M27").Value = -201733.2
$ws.Range("N27").Value = -3947.6667
$ws.Range("H45").Value = 30968
$ws.Range("I45").Value = 3980
$ws.Range("J45").Value = 39964
$ws.Range("K45").Value = 3980
$ws.Range("L45").Value = 39964
$ws.Range("M45").Value = -3573
$ws.Range("N45").Value = -40778
$ws.Range("H46").Value = 4000
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 3750
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 3750
$ws.Range("M46").Value = -4812
$ws.Range("N46").Value = -4126
$ws.Range("H55").Value = 423.4
$ws.Range("I55").Value = 380.2
$ws.Range("J55").Value = 445
$ws.Range("K55").Value = 380.2
$ws.Range("L55").Value = 445
$ws.Range("M55").Value = -207.2
$ws.Range("N55").Value = -791
$ws.Range("H123").Value = 29205.7
$ws.Range("J123").Value = 29205.7
$ws.Range("L123").Value = 29205.7
$ws.Range("N123").Value = -39005.7
$ws.Range("H136").Value = 2266.4849
$ws.Range("I136").Value = 1249.7858
$ws.Range("J136").Value = 7960
$ws.Range("K136").Value = 3749.3574
$ws.Range("L136").Value = 23880
$ws.Range("M136").Value = -1199.3574
$ws.Range("N136").Value = -28980

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 172451580
$ws.Range("I96").Value = 252650000
$ws.Range("J96").Value = 12054762
$ws.Range("K96").Value = 252650000
$ws.Range("L96").Value = 12054762
$ws.Range("M96").Value = -252648627
$ws.Range("N96").Value = -12057508
$ws.Range("H107").Value = 790.1111
$ws.Range("I107").Value = 658.2857
$ws.Range("J107").Value = 1251.5
$ws.Range("K107").Value = 1974.8571
$ws.Range("L107").Value = 3754.5
$ws.Range("M107").Value = -54.85710000000017
$ws.Range("N107").Value = -7594.5
$ws.Range("H122").Value = 4532.636
$ws.Range("I122").Value = 3377.9
$ws.Range("J122").Value = 5494.9165
$ws.Range("K122").Value = 10133.7
$ws.Range("L122").Value = 16484.7495
$ws.Range("M122").Value = -7683.700000000001
$ws.Range("N122").Value = -21384.7495
$ws.Range("H132").Value = 5377973.5
$ws.Range("I132").Value = 1104.3673
$ws.Range("J132").Value = 25644634
$ws.Range("K132").Value = 3313.1019
$ws.Range("L132").Value = 76933902
$ws.Range("M132").Value = -783.1018999999997
$ws.Range("N132").Value = -76938962
